# TestingLog.xlsx update — user, admin, doc and static controller testing results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controller - testing")
$ws.Activate()

# Excel color ints (0x00BBGGRR) used by the existing fills in this workbook
$colorYellow = 65535     # FFFF00
$colorGreen  = 5287936   # 00B050
$colorRed    = 255       # FF0000

# -----------------------------------------------------------------
# Set the brand-new text values first, in the exact order they were
# originally authored, so the shared-string table comes out the same.
# -----------------------------------------------------------------
$ws.Range("E9").Value  = ":update remoces password, when it is sent into update, does not properly redirect to :edit"
$ws.Range("E10").Value = ":delete method not found"
$ws.Range("E11").Value = ":create does not increment or redirect"
$ws.Range("E3").Value  = "same as admin_controller, except also has undefined local variable or method patient_params"
$ws.Range("E6").Value  = "destroy method is missing, update method does not work, and the:show method does render the appropriate template"
$ws.Range("E2").Value  = "i dont know how to test this, what exactly is it doing in a backend context so i can check its occuring"
$ws.Range("B2").Value  = "no-read comment"
$ws.Range("E5").Value  = "need to test content of tthe views"

# -----------------------------------------------------------------
# Remaining cell values (these reuse already-existing shared strings)
# -----------------------------------------------------------------
$ws.Range("B3").Value = "yes"
$ws.Range("D3").Value = "no"
$ws.Range("B5").Value = "yes"
$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = "yes"
$ws.Range("E7").Value = "you are missing the index view, so at this point in time all test which require an index fail. Good news this controller has a destroy method"
$ws.Range("D8").Value = ":show"
$ws.Range("E8").Value = ":index does not work (returns nil), :new does not work returns nil, :edit doesnt work (does not redirect and returns nil)"
$ws.Range("F8").ClearContents()

# -----------------------------------------------------------------
# Cell formatting: fills + vertical alignment + wrap text
# -----------------------------------------------------------------
$ws.Range("A2").Interior.Color = $colorRed
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E2").WrapText = $true

$ws.Range("A3").Interior.Color = $colorRed
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").WrapText = $true

$ws.Range("A5").Interior.Color = $colorGreen
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("C5").VerticalAlignment = -4108
$ws.Range("D5").VerticalAlignment = -4108
$ws.Range("E5").VerticalAlignment = -4108
$ws.Range("E5").WrapText = $true

$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# -----------------------------------------------------------------
# Row heights
# -----------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30

# -----------------------------------------------------------------
# Selection / view state
# -----------------------------------------------------------------
[void]$ws.Range("C6").Select()
